$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3200
$ws.Range("J40").Value = 3200
$ws.Range("L40").Value = 3200
$ws.Range("N40").Value = -3550
$ws.Range("H111").Value = 2114.5186
$ws.Range("I111").Value = 816.8
$ws.Range("J111").Value = 5822.2856
$ws.Range("K111").Value = 2450.4
$ws.Range("L111").Value = 17466.8568
$ws.Range("M111").Value = 616.6000000000004
$ws.Range("N111").Value = -23600.8568
$ws.Range("H112").Value = 4313337
$ws.Range("J112").Value = 3099.2856
$ws.Range("L112").Value = 9297.856800000001
$ws.Range("N112").Value = -11513.8568
$ws.Range("H113").Value = 26320518
$ws.Range("J113").Value = 5121.5
$ws.Range("L113").Value = 5121.5
$ws.Range("N113").Value = -11629.5
$ws.Range("H116").Value = 4710.1816
$ws.Range("H137").Value = 18859.355
$ws.Range("I137").Value = 2106.6765
$ws.Range("K137").Value = 6320.029500000001
$ws.Range("M137").Value = -3770.029500000001
$ws.Range("H138").Value = 2813.9812
$ws.Range("I138").Value = 2472.4546
$ws.Range("J138").Value = 2903.4285
$ws.Range("K138").Value = 7417.3638
$ws.Range("L138").Value = 8710.2855
$ws.Range("M138").Value = -2277.3638
$ws.Range("N138").Value = -18990.2855

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1617.5
$ws.Range("I2").Value = 1070.1666
$ws.Range("J2").Value = 2602.7
$ws.Range("K2").Value = 1070.1666
$ws.Range("L2").Value = 2602.7
$ws.Range("M2").Value = -957.1666
$ws.Range("N2").Value = -2828.7
$ws.Range("H32").Value = 18756.885
$ws.Range("I32").Value = 19609.928
$ws.Range("K32").Value = 19609.928
$ws.Range("M32").Value = -19322.928
$ws.Range("H45").Value = 3089.2068
$ws.Range("I45").Value = 2200.9
$ws.Range("J45").Value = 3556.7368
$ws.Range("K45").Value = 2200.9
$ws.Range("L45").Value = 3556.7368
$ws.Range("M45").Value = -1823.9
$ws.Range("N45").Value = -4310.736800000001
$ws.Range("H61").Value = 754744.25
$ws.Range("I61").Value = 1640821.4
$ws.Range("J61").Value = 4986.6924
$ws.Range("K61").Value = 1640821.4
$ws.Range("L61").Value = 4986.6924
$ws.Range("M61").Value = -1640609.4
$ws.Range("N61").Value = -5410.6924
$ws.Range("H74").Value = 3115.238
$ws.Range("J74").Value = 1915.5
$ws.Range("L74").Value = 1915.5
$ws.Range("N74").Value = -3663.5
$ws.Range("H77").Value = 3115.238
$ws.Range("J77").Value = 1915.5
$ws.Range("L77").Value = 9577.5
$ws.Range("N77").Value = -18313.5
$ws.Range("H116").Value = 1617.5
$ws.Range("I116").Value = 1070.1666
$ws.Range("J116").Value = 2602.7
$ws.Range("K116").Value = 1070.1666
$ws.Range("L116").Value = 2602.7
$ws.Range("M116").Value = 1223.8334
$ws.Range("N116").Value = -7190.7
$ws.Range("H132").Value = 20584.555
$ws.Range("I132").Value = 1854.0714
$ws.Range("K132").Value = 5562.2142
$ws.Range("M132").Value = -3032.2142
$ws.Range("H136").Value = 754744.25
$ws.Range("I136").Value = 1640821.4
$ws.Range("J136").Value = 4986.6924
$ws.Range("K136").Value = 4922464.199999999
$ws.Range("L136").Value = 14960.0772
$ws.Range("M136").Value = -4919914.199999999
$ws.Range("N136").Value = -20060.0772

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1617.5
$ws.Range("I3").Value = 1070.1666
$ws.Range("J3").Value = 2602.7
$ws.Range("K3").Value = 1070.1666
$ws.Range("L3").Value = 2602.7
$ws.Range("M3").Value = -956.1666
$ws.Range("N3").Value = -2830.7
$ws.Range("H134").Value = 49181.59
$ws.Range("I134").Value = 59833.055
$ws.Range("J134").Value = 1250
$ws.Range("K134").Value = 179499.165
$ws.Range("L134").Value = 3750
$ws.Range("M134").Value = -176964.165
$ws.Range("N134").Value = -8820

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10608.021
$ws.Range("I31").Value = 16410
$ws.Range("J31").Value = 3424.6191
$ws.Range("K31").Value = 16410
$ws.Range("L31").Value = 3424.6191
$ws.Range("M31").Value = -16115
$ws.Range("N31").Value = -4014.6191
$ws.Range("H34").Value = 10608.021
$ws.Range("I34").Value = 16410
$ws.Range("J34").Value = 3424.6191
$ws.Range("K34").Value = 16410
$ws.Range("L34").Value = 3424.6191
$ws.Range("M34").Value = -16208
$ws.Range("N34").Value = -3828.6191
$ws.Range("H58").Value = 23018.912
$ws.Range("I58").Value = 1309
$ws.Range("J58").Value = 126141
$ws.Range("K58").Value = 1309
$ws.Range("L58").Value = 126141
$ws.Range("M58").Value = -1106
$ws.Range("N58").Value = -126547
$ws.Range("H132").Value = 45183.4
$ws.Range("I132").Value = 70036.11
$ws.Range("K132").Value = 210108.33
$ws.Range("M132").Value = -207578.33
$ws.Range("H134").Value = 3939.5625
$ws.Range("I134").Value = 698.8276
$ws.Range("J134").Value = 35266.668
$ws.Range("K134").Value = 2096.4828
$ws.Range("L134").Value = 105800.004
$ws.Range("M134").Value = 438.5172000000002
$ws.Range("N134").Value = -110870.004
$ws.Range("H136").Value = 23018.912
$ws.Range("I136").Value = 1309
$ws.Range("J136").Value = 126141
$ws.Range("K136").Value = 3927
$ws.Range("L136").Value = 378423
$ws.Range("M136").Value = -1377
$ws.Range("N136").Value = -383523

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 8471
$ws.Range("I68").Value = 949.6667
$ws.Range("J68").Value = 10522.272
$ws.Range("K68").Value = 2849.0001
$ws.Range("L68").Value = 31566.816
$ws.Range("M68").Value = -2038.0001
$ws.Range("N68").Value = -33188.81600000001
$ws.Range("H71").Value = 8471
$ws.Range("I71").Value = 949.6667
$ws.Range("J71").Value = 10522.272
$ws.Range("K71").Value = 8547.0003
$ws.Range("L71").Value = 94700.448
$ws.Range("M71").Value = -4491.0003
$ws.Range("N71").Value = -102812.448
$ws.Range("H97").Value = 858.4545000000001
$ws.Range("J97").Value = 1215.5
$ws.Range("L97").Value = 3646.5
$ws.Range("N97").Value = -4638.5
$ws.Range("H104").Value = 10000
$ws.Range("J104").Value = 10000
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -35242
$ws.Range("H113").Value = 5003.875
$ws.Range("I113").Value = 6306.353
$ws.Range("J113").Value = 1840.7142
$ws.Range("K113").Value = 18919.059
$ws.Range("L113").Value = 5522.142599999999
$ws.Range("M113").Value = -16749.059
$ws.Range("N113").Value = -9862.142599999999
$ws.Range("H131").Value = 814.5599999999999
$ws.Range("I131").Value = 600
$ws.Range("J131").Value = 825.85266
$ws.Range("K131").Value = 1800
$ws.Range("L131").Value = 2477.55798
$ws.Range("M131").Value = 3240
$ws.Range("N131").Value = -12557.55798
$ws.Range("H140").Value = 1414.9546
$ws.Range("I140").Value = 1186.7
$ws.Range("J140").Value = 3697.5
$ws.Range("K140").Value = 3560.1
$ws.Range("L140").Value = 11092.5
$ws.Range("M140").Value = 1619.9
$ws.Range("N140").Value = -21452.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 20300
$ws.Range("J15").Value = 20300
$ws.Range("L15").Value = 20300
$ws.Range("N15").Value = -20876
$ws.Range("H81").Value = 20300
$ws.Range("J81").Value = 20300
$ws.Range("L81").Value = 20300
$ws.Range("N81").Value = -22296
$ws.Range("H84").Value = 20300
$ws.Range("J84").Value = 20300
$ws.Range("L84").Value = 60900
$ws.Range("N84").Value = -70884
$ws.Range("H122").Value = 2527.5264
$ws.Range("I122").Value = 2344.6875
$ws.Range("J122").Value = 3502.6667
$ws.Range("K122").Value = 7034.0625
$ws.Range("L122").Value = 10508.0001
$ws.Range("M122").Value = -4584.0625
$ws.Range("N122").Value = -15408.0001
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 72179.59
$ws.Range("I132").Value = 81381.84
$ws.Range("J132").Value = 58887.445
$ws.Range("K132").Value = 244145.52
$ws.Range("L132").Value = 176662.335
$ws.Range("M132").Value = -241615.52
$ws.Range("N132").Value = -181722.335
$ws.Range("H135").Value = 51585
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7100.643
$ws.Range("I61").Value = 3486.2856
$ws.Range("K61").Value = 3486.2856
$ws.Range("M61").Value = -3284.2856
$ws.Range("H110").Value = 2031379.8
$ws.Range("J110").Value = 2031379.8
$ws.Range("L110").Value = 2031379.8
$ws.Range("N110").Value = -2039559.8
$ws.Range("H113").Value = 7100.643
$ws.Range("I113").Value = 3486.2856
$ws.Range("K113").Value = 3486.2856
$ws.Range("M113").Value = -1316.2856
$ws.Range("H132").Value = 1577.3889
$ws.Range("I132").Value = 1075.9333
$ws.Range("K132").Value = 3227.7999
$ws.Range("M132").Value = -697.7999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1446.3529
$ws.Range("I126").Value = 1214.5385
$ws.Range("J126").Value = 2199.75
$ws.Range("K126").Value = 3643.6155
$ws.Range("L126").Value = 6599.25
$ws.Range("M126").Value = -1173.6155
$ws.Range("N126").Value = -11539.25
$ws.Range("H132").Value = 2217.5366
$ws.Range("I132").Value = 2046.2142
$ws.Range("J132").Value = 2586.5386
$ws.Range("K132").Value = 6138.642599999999
$ws.Range("L132").Value = 7759.6158
$ws.Range("M132").Value = -3608.642599999999
$ws.Range("N132").Value = -12819.6158
$ws.Range("H136").Value = 1078.5555
$ws.Range("I136").Value = 674.86957
$ws.Range("J136").Value = 3399.75
$ws.Range("K136").Value = 2024.60871
$ws.Range("L136").Value = 10199.25
$ws.Range("M136").Value = 525.39129
$ws.Range("N136").Value = -15299.25
